# Update "want to go" counts (column F) on the 展览 (Exhibition), 演出 (Performance)
# and 全部类型 (All types) sheets to match the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 20111
$ws1.Range("F9").Value = 7583
$ws1.Range("F38").Value = 12649
$ws1.Range("F43").Value = 270
$ws1.Range("F44").Value = 370

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 182

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 20112
$ws4.Range("F9").Value = 7583
$ws4.Range("F31").Value = 182
$ws4.Range("F40").Value = 12649
$ws4.Range("F45").Value = 270
$ws4.Range("F46").Value = 370
